$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.854.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.545.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.18%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.543.64"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.15%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.47%  "
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").Value = "  +4.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.155.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.42%  "
$ws.Range("E14").Value = "  +3.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.553.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.56%  "
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.742.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.06%  "
$ws.Range("E20").Value = "  +7.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.16%  "
$ws.Range("E23").Value = "  +6.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.692.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000116"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.38%  "
$ws.Range("E28").Value = "  +8.15%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  +5.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.556.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +23.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.20%  "
$ws.Range("E36").Value = "  +3.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "170.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("E38").Value = "  +5.21%  "
$ws.Range("E39").Value = "  +7.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0804"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.77%  "
$ws.Range("E42").Value = "  +4.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +22.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  +5.20%  "
$ws.Range("E47").Value = "  +9.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.88%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.95%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.438.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +16.75%  "
